# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#    columns get populated for both language sheets (zh-cn, de-de)
#  - The "Latest Target File" cell becomes a hyperlink back to the source doc
#  - A few columns are widened so the new, longer text fits

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c4c2392910a61cafe455de7e6d53088d2e5d758/e2e/"

$doc1 = "a745513b-cad4-409b-a5fc-fbc964c51b7f"
$doc2 = "c8113dad-defc-4e54-ad36-d9bf360d8c20"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: widen the per-language status columns (E, F) so the new,
# longer status text is fully visible.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet (row 2 = doc1, row 3 = doc2)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("J2").Value = $doc1 + ".c0b13cc25c27cf31556a90bc4baaeb0864953ff3.zh-cn.xlf"
$zhcn.Range("J3").Value = $doc2 + ".a9e36cd029b0e4fd58ae264dee550c6bc481f9f5.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-21 22:40:31"
$zhcn.Range("K3").Value = "2016-08-21 22:40:31"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $githubBase + $doc1 + ".md", "", "", $doc1 + ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $githubBase + $doc2 + ".md", "", "", $doc2 + ".md")

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet (row 2 = doc1, row 3 = doc2)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = $doc1 + ".c0b13cc25c27cf31556a90bc4baaeb0864953ff3.de-de.xlf"
$dede.Range("J3").Value = $doc2 + ".a9e36cd029b0e4fd58ae264dee550c6bc481f9f5.de-de.xlf"

$dede.Range("K2").Value = "2016-08-21 22:40:37"
$dede.Range("K3").Value = "2016-08-21 22:40:37"

$dede.Hyperlinks.Add($dede.Range("I2"), $githubBase + $doc1 + ".md", "", "", $doc1 + ".md")
$dede.Hyperlinks.Add($dede.Range("I3"), $githubBase + $doc2 + ".md", "", "", $doc2 + ".md")

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
